$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 325; existing rows 325:354 shift down to 326:355
$ws.Rows(325).Insert()

# Populate the newly inserted row 325 with its data
$ws.Range("A325").Value = 4
$ws.Range("B325").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C325").Value = "Los Lagos"
$ws.Range("D325").Value = 44585
$ws.Range("D325").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E325").Value = 10
$ws.Range("F325").Value = 100112006
$ws.Range("G325").Value = "Repollo"
$ws.Range("H325").Value = "Crespo record"
$ws.Range("I325").Value = "Primera"
$ws.Range("J325").Value = 500
$ws.Range("K325").Value = 1500
$ws.Range("L325").Value = 1500
$ws.Range("M325").Value = 1500
$ws.Range("N325").Value = "$/unidad"
$ws.Range("O325").Value = "Región Metropolitana"
$ws.Range("P325").Value = 1500
$ws.Range("Q325").Value = 1
$ws.Range("R325").Value = "Hortaliza"
